$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1403
$ws.Range("F4").Value = 26049
$ws.Range("F5").Value = 573
$ws.Range("F6").Value = 244
$ws.Range("F7").Value = 583
$ws.Range("F10").Value = 231
$ws.Range("F11").Value = 346
$ws.Range("F12").Value = 209
$ws.Range("F13").Value = 179
$ws.Range("F15").Value = 279
$ws.Range("F16").Value = 33
$ws.Range("F17").Value = 348
$ws.Range("F18").Value = 48
$ws.Range("F19").Value = 1489
$ws.Range("F20").Value = 167
$ws.Range("F21").Value = 12
$ws.Range("F23").Value = 92
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 220
$ws.Range("F4").Value = 3
$ws.Range("F6").Value = 114
$ws.Range("F8").Value = 104
$ws.Range("F9").Value = 104
$ws.Range("F10").Value = 426
$ws.Range("F13").Value = 10
$ws.Range("F15").Value = 41
$ws.Range("F17").Value = 26
$ws.Range("F18").Value = 1
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4927
$ws.Range("F3").Value = 191
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1403
$ws.Range("F4").Value = 4927
$ws.Range("F5").Value = 191
$ws.Range("F6").Value = 26049
$ws.Range("F7").Value = 573
$ws.Range("F9").Value = 244
$ws.Range("F10").Value = 220
$ws.Range("F11").Value = 583
$ws.Range("F12").Value = 3
$ws.Range("F15").Value = 115
$ws.Range("F16").Value = 115
$ws.Range("F18").Value = 104
$ws.Range("F19").Value = 104
$ws.Range("F20").Value = 426
$ws.Range("F23").Value = 231
$ws.Range("F24").Value = 346
$ws.Range("F25").Value = 209
$ws.Range("F26").Value = 179
$ws.Range("F29").Value = 279
$ws.Range("F30").Value = 33
$ws.Range("F31").Value = 10
$ws.Range("F33").Value = 348
$ws.Range("F34").Value = 48
$ws.Range("F35").Value = 41
$ws.Range("F36").Value = 1489
$ws.Range("F37").Value = 167
$ws.Range("F39").Value = 12
$ws.Range("F41").Value = 92
$ws.Range("F43").Value = 26
$ws.Range("F44").Value = 1
